$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.297.91"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "1.790.81"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'226.46"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  +1.20%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'32.71"
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("E9").Value = "  +0.68%  "
$ws.Range("D10").Value = "'0.0688"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "2.048.57"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").Value = "1.784.88"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "'0.631"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").Value = "34.305.74"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").Value = "'4.26"
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").Value = "'68.29"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").Value = "0.0₃0792"
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "'243.82"
$ws.Range("D21").Value = "'11.28"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").Value = "'165.40"
$ws.Range("E25").Value = "  +2.72%  "
$ws.Range("D26").Value = "'7.28"
$ws.Range("E26").Value = "  +1.93%  "
$ws.Range("D27").Value = "'16.44"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("E30").Value = "  +5.82%  "
$ws.Range("D31").Value = "'0.0522"
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.23"
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'3.78"
$ws.Range("E33").Value = "  +0.67%  "
$ws.Range("E34").Value = "  +0.69%  "
$ws.Range("E35").Value = "  -1.31%  "
$ws.Range("D36").Value = "1.398.63"
$ws.Range("E36").Value = "  -3.15%  "
$ws.Range("D37").Value = "'0.667"
$ws.Range("E37").Value = "  +1.64%  "
$ws.Range("E38").Value = "  +1.64%  "
$ws.Range("E39").Value = "  -0.63%  "
$ws.Range("D40").Value = "'84.85"
$ws.Range("E40").Value = "  +3.58%  "
$ws.Range("E41").Value = "  +4.06%  "
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("D43").Value = "'0.933"
$ws.Range("E43").Value = "  +1.59%  "
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("E45").Value = "  +1.75%  "
$ws.Range("E46").Value = "  +2.92%  "
$ws.Range("D47").Value = "'5.99"
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("D48").Value = "1.950.14"
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("D49").Value = "'104.65"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("E51").Value = "  -1.30%  "
